# Fruta / hortaliza, semanal
# Insert a new weekly record at row 151 of "Sheet1" for the
# "Vega Modelo de Temuco - Espinaca" sheet, pushing the existing
# rows 151-156 down to 152-157.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 151 (shifts 151..156 -> 152..157)
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with this week's data
$ws.Cells.Item(151, 1).Value2 = 10
$ws.Cells.Item(151, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(151, 3).Value2 = "La Araucanía"
$ws.Cells.Item(151, 4).Value2 = 44753
$ws.Cells.Item(151, 5).Value2 = 9
$ws.Cells.Item(151, 6).Value2 = 100112012
$ws.Cells.Item(151, 7).Value2 = "Espinaca"
$ws.Cells.Item(151, 8).Value2 = "Sin especificar"
$ws.Cells.Item(151, 9).Value2 = "Segunda"
$ws.Cells.Item(151, 10).Value2 = 30
$ws.Cells.Item(151, 11).Value2 = 10000
$ws.Cells.Item(151, 12).Value2 = 10000
$ws.Cells.Item(151, 13).Value2 = 10000
$ws.Cells.Item(151, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(151, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(151, 16).Value2 = 3333
$ws.Cells.Item(151, 17).Value2 = 3
$ws.Cells.Item(151, 18).Value2 = "Hortaliza"
